# Actualizar ligas y agregar Liga Argentina actualizada
# Appends new fixtures (rows 130-137) to Sheet1 of the Primera B 2025 workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=130; A="2025-07-19"; B="Cobreloa"; C="San Marcos de Arica"; D=2; E=0; F=1348361; G=5; H=4; I=3; J=6; K=0; L=1; M=1; N=0; O=1; P=0; Q=50; R=50; S="L" },
    @{ Row=131; A="2025-07-19"; B="Deportes Santa Cruz"; C="Recoleta"; D=0; E=0; F=1348366; G=9; H=9; I=4; J=1; K=0; L=0; M=0; N=0; O=0; P=0; Q=47; R=53; S="E" },
    @{ Row=132; A="2025-07-19"; B="Concepción"; C="Antofagasta"; D=0; E=0; F=1348367; G=12; H=0; I=3; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=57; R=43; S="E" },
    @{ Row=133; A="2025-07-20"; B="Deportes Copiapo"; C="Curico Unido"; D=1; E=1; F=1348362; G=8; H=4; I=3; J=5; K=0; L=0; M=1; N=0; O=0; P=1; Q=55; R=45; S="E" },
    @{ Row=134; A="2025-07-20"; B="Magallanes"; C="Universidad de Concepcion"; D=1; E=2; F=1348365; G=7; H=2; I=2; J=1; K=0; L=0; M=1; N=0; O=0; P=2; Q=52; R=48; S="V" },
    @{ Row=135; A="2025-07-20"; B="Deportes Temuco"; C="Santiago Morning"; D=2; E=1; F=1348368; G=3; H=9; I=3; J=6; K=0; L=1; M=0; N=0; O=2; P=1; Q=55; R=45; S="L" },
    @{ Row=136; A="2025-07-21"; B="San Luis"; C="Rangers de Talca"; D=1; E=2; F=1348363; G=4; H=3; I=2; J=2; K=1; L=0; M=0; N=1; O=1; P=1; Q=56; R=44; S="V" },
    @{ Row=137; A="2025-07-22"; B="Santiago Wanderers"; C="Union San Felipe"; D=1; E=3; F=1348364; G=3; H=3; I=7; J=3; K=1; L=0; M=1; N=0; O=0; P=3; Q=45; R=55; S="V" }
)

$numericCols = @{ D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18 }

foreach ($r in $rows) {
    # Column A holds a date-formatted text string ("yyyy-mm-dd"), matching the
    # rest of the sheet where every date is stored as plain text, never a
    # real date serial. A plain `.Value = "2025-07-19"` assignment gets
    # auto-parsed into a date by Excel's input heuristics, so instead we
    # write it as a text formula (forcing a string result), then demote the
    # formula to a static value via copy / paste-values (xlPasteValues),
    # which carries the string through verbatim without re-triggering the
    # date auto-detection that a direct `.Value = "<text>"` assignment would.
    $cellA = $ws.Cells.Item($r.Row, 1)
    $cellA.Formula = '="' + $r.A + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C

    foreach ($col in $numericCols.Keys) {
        $ws.Cells.Item($r.Row, $numericCols[$col]).Value = $r[$col]
    }

    $ws.Cells.Item($r.Row, 19).Value = $r.S
}

# Match the author's final view/selection state.
$excel.ActiveWindow.ScrollRow = 124
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G136").Select()
